$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.16%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.95%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.334"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.56%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07427"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "11.22%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.779"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.57%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.691"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "8.25%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.566"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "15.32%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9126"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01667"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2,481.15%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07609"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "13.97%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08075"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.72%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03014"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.73%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09851"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "9.55%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001524"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.06%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04554"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.90%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006354"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.25%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.497"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.44%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.237"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.77%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.59%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1335"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.94%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.78%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1627"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.80%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001215"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.19%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004498"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.99%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.44%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001739"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "7.51%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04509"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "7.09%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007201"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.71%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.78%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002259"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "14.10%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01371"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.86%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006111"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.00%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.892"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.02%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01299"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.62%"
